# Auto-generated edit script applying numeric value updates per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1422.0217
$ws.Range("J17").Value = 1439.6
$ws.Range("L17").Value = 4318.799999999999
$ws.Range("N17").Value = -4654.799999999999
$ws.Range("H33").Value = 761.9
$ws.Range("I33").Value = 183.9
$ws.Range("K33").Value = 183.9
$ws.Range("M33").Value = 45.09999999999999
$ws.Range("H58").Value = 20838330
$ws.Range("I58").Value = 25000514
$ws.Range("K58").Value = 75001542
$ws.Range("M58").Value = -75001392
$ws.Range("H62").Value = 35671
$ws.Range("I62").Value = 164599.4
$ws.Range("K62").Value = 164599.4
$ws.Range("M62").Value = -163975.4
$ws.Range("H65").Value = 35671
$ws.Range("I65").Value = 164599.4
$ws.Range("K65").Value = 822997
$ws.Range("M65").Value = -819877
$ws.Range("H86").Value = 3045.5334
$ws.Range("J86").Value = 3138.6667
$ws.Range("L86").Value = 3138.6667
$ws.Range("N86").Value = -5384.6667
$ws.Range("H88").Value = 3289.8635
$ws.Range("I88").Value = 699.6667
$ws.Range("J88").Value = 3698.842
$ws.Range("K88").Value = 699.6667
$ws.Range("L88").Value = 3698.842
$ws.Range("M88").Value = -293.6667
$ws.Range("N88").Value = -4510.842000000001
$ws.Range("H89").Value = 3045.5334
$ws.Range("J89").Value = 3138.6667
$ws.Range("L89").Value = 15693.3335
$ws.Range("N89").Value = -26925.3335
$ws.Range("H91").Value = 3289.8635
$ws.Range("I91").Value = 699.6667
$ws.Range("J91").Value = 3698.842
$ws.Range("K91").Value = 699.6667
$ws.Range("L91").Value = 3698.842
$ws.Range("M91").Value = 704.3333
$ws.Range("N91").Value = -6506.842000000001
$ws.Range("H141").Value = 3218
$ws.Range("I141").Value = 2897.5
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 8692.5
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = -3512.5
$ws.Range("N141").Value = -23860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 11216.714
$ws.Range("J43").Value = 11419.333
$ws.Range("L43").Value = 11419.333
$ws.Range("N43").Value = -12045.333
$ws.Range("H63").Value = 5721.5557
$ws.Range("I63").Value = 2298.8
$ws.Range("K63").Value = 2298.8
$ws.Range("M63").Value = -1612.8
$ws.Range("H66").Value = 5721.5557
$ws.Range("I66").Value = 2298.8
$ws.Range("K66").Value = 11494
$ws.Range("M66").Value = -8062
$ws.Range("H113").Value = 46000
$ws.Range("J113").Value = 46000
$ws.Range("L113").Value = 46000
$ws.Range("N113").Value = -54678

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2991.2
$ws.Range("I86").Value = 2991.2
$ws.Range("K86").Value = 2991.2
$ws.Range("M86").Value = -1868.2
$ws.Range("H89").Value = 2991.2
$ws.Range("I89").Value = 2991.2
$ws.Range("K89").Value = 14956
$ws.Range("M89").Value = -9340
$ws.Range("H94").Value = 3393.6
$ws.Range("I94").Value = 3167.2
$ws.Range("J94").Value = 4072.8
$ws.Range("K94").Value = 3167.2
$ws.Range("L94").Value = 4072.8
$ws.Range("M94").Value = -2716.2
$ws.Range("N94").Value = -4974.8
$ws.Range("H105").Value = 2068.4211
$ws.Range("I105").Value = 2058.4
$ws.Range("J105").Value = 2106
$ws.Range("K105").Value = 2058.4
$ws.Range("L105").Value = 2106
$ws.Range("M105").Value = -311.4000000000001
$ws.Range("N105").Value = -5600
$ws.Range("H107").Value = 3157.4243
$ws.Range("I107").Value = 1213
$ws.Range("K107").Value = 1213
$ws.Range("M107").Value = 707
$ws.Range("H134").Value = 2205.2727
$ws.Range("I134").Value = 2045.32
$ws.Range("K134").Value = 6135.96
$ws.Range("M134").Value = -3600.96
$ws.Range("H140").Value = 61750
$ws.Range("J140").Value = 61750
$ws.Range("L140").Value = 61750
$ws.Range("N140").Value = -72110

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4781.184
$ws.Range("I31").Value = 2454.96
$ws.Range("J31").Value = 9254.691999999999
$ws.Range("K31").Value = 2454.96
$ws.Range("L31").Value = 9254.691999999999
$ws.Range("M31").Value = -2159.96
$ws.Range("N31").Value = -9844.691999999999
$ws.Range("H34").Value = 4781.184
$ws.Range("I34").Value = 2454.96
$ws.Range("J34").Value = 9254.691999999999
$ws.Range("K34").Value = 2454.96
$ws.Range("L34").Value = 9254.691999999999
$ws.Range("M34").Value = -2252.96
$ws.Range("N34").Value = -9658.691999999999
$ws.Range("H62").Value = 11666.667
$ws.Range("I62").Value = 10750
$ws.Range("J62").Value = 13500
$ws.Range("K62").Value = 10750
$ws.Range("L62").Value = 13500
$ws.Range("M62").Value = -10126
$ws.Range("N62").Value = -14748
$ws.Range("H65").Value = 11666.667
$ws.Range("I65").Value = 10750
$ws.Range("J65").Value = 13500
$ws.Range("K65").Value = 53750
$ws.Range("L65").Value = 67500
$ws.Range("M65").Value = -50630
$ws.Range("N65").Value = -73740
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61498
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -187488
$ws.Range("H88").Value = 18030.666
$ws.Range("J88").Value = 18030.666
$ws.Range("L88").Value = 18030.666
$ws.Range("N88").Value = -18842.666
$ws.Range("H91").Value = 18030.666
$ws.Range("J91").Value = 18030.666
$ws.Range("L91").Value = 18030.666
$ws.Range("N91").Value = -20838.666
$ws.Range("H107").Value = 1127.1714
$ws.Range("I107").Value = 636.4666999999999
$ws.Range("K107").Value = 636.4666999999999
$ws.Range("M107").Value = 1283.5333
$ws.Range("H122").Value = 3524.077
$ws.Range("I122").Value = 774
$ws.Range("J122").Value = 4349.1
$ws.Range("K122").Value = 2322
$ws.Range("L122").Value = 13047.3
$ws.Range("M122").Value = 128
$ws.Range("N122").Value = -17947.3
$ws.Range("H132").Value = 2664.5715
$ws.Range("I132").Value = 2745.6562
$ws.Range("K132").Value = 8236.9686
$ws.Range("M132").Value = -5706.9686

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 611166.5
$ws.Range("I4").Value = 266140.16
$ws.Range("K4").Value = 798420.48
$ws.Range("M4").Value = -798308.48
$ws.Range("H23").Value = 4430.227
$ws.Range("I23").Value = 4414.6665
$ws.Range("J23").Value = 4448.9
$ws.Range("K23").Value = 13243.9995
$ws.Range("L23").Value = 13346.7
$ws.Range("M23").Value = -13008.9995
$ws.Range("N23").Value = -13816.7
$ws.Range("H38").Value = 634.2222
$ws.Range("J38").Value = 938.1667
$ws.Range("L38").Value = 2814.5001
$ws.Range("N38").Value = -3508.5001
$ws.Range("H39").Value = 1999.3334
$ws.Range("I39").Value = 1999.5
$ws.Range("J39").Value = 1999
$ws.Range("K39").Value = 5998.5
$ws.Range("L39").Value = 5997
$ws.Range("M39").Value = -5704.5
$ws.Range("N39").Value = -6585
$ws.Range("H140").Value = 1990
$ws.Range("I140").Value = 1990
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 5970
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -790

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 939.61536
$ws.Range("I107").Value = 792.4
$ws.Range("J107").Value = 1430.3334
$ws.Range("K107").Value = 792.4
$ws.Range("L107").Value = 1430.3334
$ws.Range("M107").Value = 1127.6
$ws.Range("N107").Value = -5270.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 3808.3333
$ws.Range("J22").Value = 3240.8333
$ws.Range("K22").Value = 3808.3333
$ws.Range("L22").Value = 3240.8333
$ws.Range("M22").Value = -3513.3333
$ws.Range("N22").Value = -3830.8333
$ws.Range("I27").Value = 3808.3333
$ws.Range("J27").Value = 3240.8333
$ws.Range("K27").Value = 3808.3333
$ws.Range("L27").Value = 3240.8333
$ws.Range("M27").Value = -3701.3333
$ws.Range("N27").Value = -3454.8333
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H68").Value = 8374.8125
$ws.Range("I68").Value = 6249.75
$ws.Range("J68").Value = 9083.166999999999
$ws.Range("K68").Value = 6249.75
$ws.Range("L68").Value = 9083.166999999999
$ws.Range("M68").Value = -5500.75
$ws.Range("N68").Value = -10581.167
$ws.Range("H71").Value = 8374.8125
$ws.Range("I71").Value = 6249.75
$ws.Range("J71").Value = 9083.166999999999
$ws.Range("K71").Value = 31248.75
$ws.Range("L71").Value = 45415.835
$ws.Range("M71").Value = -27504.75
$ws.Range("N71").Value = -52903.835
$ws.Range("H132").Value = 4099.5835
$ws.Range("I132").Value = 4699.1665
$ws.Range("K132").Value = 14097.4995
$ws.Range("M132").Value = -11567.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 31222
$ws.Range("J49").Value = 34000
$ws.Range("L49").Value = 34000
$ws.Range("N49").Value = -34460
$ws.Range("H74").Value = 9009.556
$ws.Range("J74").Value = 9009.556
$ws.Range("L74").Value = 9009.556
$ws.Range("N74").Value = -10881.556
$ws.Range("H77").Value = 9009.556
$ws.Range("J77").Value = 9009.556
$ws.Range("L77").Value = 27028.668
$ws.Range("N77").Value = -36388.66800000001
$ws.Range("H81").Value = 2916.5
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 2916.5
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
$ws.Range("H122").Value = 5884.3687
$ws.Range("I122").Value = 2100.375
$ws.Range("K122").Value = 6301.125
$ws.Range("M122").Value = -3851.125

# ---- Clear cells that were removed in the target revision ----
$wb.Worksheets.Item("CUL").Range("N140").ClearContents() | Out-Null
$wb.Worksheets.Item("LTW").Range("N38").ClearContents() | Out-Null
